$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.486.92"
$ws.Range("E2").Value = "  -1.15%  "
$ws.Range("D3").Value = "1.911.59"
$ws.Range("E3").Value = "  -1.52%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'238.90"
$ws.Range("E5").Value = "  -1.74%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "'0.4765"
$ws.Range("E7").Value = "  -2.69%  "
$ws.Range("D8").Value = "'0.2835"
$ws.Range("E8").Value = "  -3.30%  "
$ws.Range("D9").Value = "'0.06700"
$ws.Range("E9").Value = "  -2.90%  "
$ws.Range("D10").Value = "'19.45"
$ws.Range("E10").Value = "  +1.03%  "
$ws.Range("D11").Value = "'103.29"
$ws.Range("E11").Value = "  -2.18%  "
$ws.Range("D12").Value = "'0.07756"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "1.914.15"
$ws.Range("E13").Value = "  -1.44%  "
$ws.Range("D14").Value = "'5.186"
$ws.Range("E14").Value = "  -3.39%  "
$ws.Range("D15").Value = "'0.6676"
$ws.Range("E15").Value = "  -4.88%  "
$ws.Range("D16").Value = "'276.77"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("D17").Value = "30.493.87"
$ws.Range("D18").Value = "'0.9992"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").Value = "'0.000007479"
$ws.Range("E19").Value = "  -3.23%  "
$ws.Range("D20").Value = "'12.62"
$ws.Range("E20").Value = "  -3.47%  "
$ws.Range("D21").Value = "'5.384"
$ws.Range("E21").Value = "  -2.51%  "
$ws.Range("D22").Value = "'0.4626"
$ws.Range("E22").Value = "  -8.20%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "'6.276"
$ws.Range("E24").Value = "  -4.26%  "
$ws.Range("D25").Value = "'9.340"
$ws.Range("E25").Value = "  -4.61%  "
$ws.Range("D26").Value = "'166.26"
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").Value = "'19.24"
$ws.Range("E27").Value = "  -1.83%  "
$ws.Range("D28").Value = "'2.069"
$ws.Range("E28").Value = "  -4.15%  "
$ws.Range("D29").Value = "'1.380"
$ws.Range("D30").Value = "'0.09984"
$ws.Range("E30").Value = "  -3.71%  "
$ws.Range("D31").Value = "'4.613"
$ws.Range("E31").Value = "  +0.71%  "
$ws.Range("D32").Value = "'1.509"
$ws.Range("E32").Value = "  -3.20%  "
$ws.Range("D33").Value = "'4.246"
$ws.Range("E33").Value = "  -3.22%  "
$ws.Range("D34").Value = "'0.04686"
$ws.Range("E34").Value = "  -4.08%  "
$ws.Range("D35").Value = "'0.7270"
$ws.Range("E35").Value = "  -3.96%  "
$ws.Range("D36").Value = "'1.113"
$ws.Range("E36").Value = "  -3.39%  "
$ws.Range("D37").Value = "'2.712"
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("D38").Value = "'0.01904"
$ws.Range("E38").Value = "  -4.78%  "
$ws.Range("D39").Value = "'2.606"
$ws.Range("E39").Value = "  -1.84%  "
$ws.Range("D40").Value = "'6.332"
$ws.Range("E40").Value = "  -2.14%  "
$ws.Range("D41").Value = "'74.22"
$ws.Range("E41").Value = "  -6.45%  "
$ws.Range("D42").Value = "'1.962"
$ws.Range("E42").Value = "  -6.02%  "
$ws.Range("D43").Value = "'0.8585"
$ws.Range("E43").Value = "  -6.26%  "
$ws.Range("D44").Value = "'105.76"
$ws.Range("E44").Value = "  -1.99%  "
$ws.Range("D45").Value = "'0.4257"
$ws.Range("E45").Value = "  -3.81%  "
$ws.Range("D46").Value = "'0.9999"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").Value = "'7.418"
$ws.Range("E47").Value = "  -3.45%  "
$ws.Range("D48").Value = "'953.46"
$ws.Range("E48").Value = "  -2.52%  "
$ws.Range("D49").Value = "'0.1207"
$ws.Range("E49").Value = "  -2.88%  "
$ws.Range("D50").Value = "'34.65"
$ws.Range("E50").Value = "  -4.00%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05796"
$ws.Range("E51").Value = "  +0.39%  "
